$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New shared strings used by the new rows (B7/D7 labels)
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "Tower"
$ws.Range("D7").Value = "Nestor pure arena"

# ---------------------------------------------------------------------------
# New data rows (8 and 9) with the distance formula, mirroring rows 3/4
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = 612.5
$ws.Range("C8").Value = 532.5
$ws.Range("D8").Value = 287
$ws.Range("E8").Value = 344
$ws.Range("F8").Formula = "=SQRT(SUMXMY2(B8:C8,D8:E8))"

$ws.Range("B9").Value = 150.5
$ws.Range("C9").Value = -146.5
$ws.Range("D9").Value = -22
$ws.Range("E9").Value = 171.25
$ws.Range("F9").Formula = "=SQRT(SUMXMY2(B9:C9,D9:E9))"

# Ratio formulas (mirrors F5/F6 for the new pair of points)
$ws.Range("F11").Formula = "=F9/F8"
$ws.Range("F12").Formula = "=F8/F9"

# ---------------------------------------------------------------------------
# Column widths (A-F) now explicit/custom
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10.666666666666666
$ws.Columns.Item(2).ColumnWidth = 33.666666666666664
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws.Columns.Item(4).ColumnWidth = 35.166666666666664
$ws.Columns.Item(5).ColumnWidth = 10.666666666666666
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668

# ---------------------------------------------------------------------------
# Number format ("#,##0.000") kept on the same cells it already decorated
# ---------------------------------------------------------------------------
$numFmtCells = "G3","H3","I3","K3","L3","G4","K4","L4","B5","C5","H5","K5","L5","H7","I7","H8","I8"
foreach ($addr in $numFmtCells) {
    $ws.Range($addr).NumberFormat = "#,##0.000"
}
$ws.Range("B3:F3").NumberFormat = "#,##0.000"
$ws.Range("B4:F4").NumberFormat = "#,##0.000"

# Header cells keep their centred alignment
$ws.Range("B1:E1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Selection moves to F12
# ---------------------------------------------------------------------------
[void]$ws.Range("F12").Select()
